$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.488.51"
$ws.Range("E2").Value = "  -2.93%  "
$ws.Range("D3").Value = "2.999.74"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "545.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "2.995.63"
$ws.Range("E8").Value = "  -2.41%  "
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.66%  "
$ws.Range("E11").Value = "  -8.47%  "
$ws.Range("E12").Value = "  -3.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000217"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.82%  "
$ws.Range("D15").Value = "3.482.28"
$ws.Range("E15").Value = "  -2.76%  "
$ws.Range("D16").Value = "61.514.08"
$ws.Range("E16").Value = "  -2.92%  "
$ws.Range("E17").Value = "  -2.82%  "
$ws.Range("D18").Value = "3.002.83"
$ws.Range("E18").Value = "  -2.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.76%  "
$ws.Range("E22").Value = "  -5.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.18%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.15%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.86%  "
$ws.Range("E32").Value = "  -4.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "54.63"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "436.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.40%  "
$ws.Range("D38").Value = "3.122.48"
$ws.Range("E38").Value = "  -4.04%  "
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("E40").Value = "  -5.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.115"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.48%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.240"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.62%  "
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("E48").Value = "  -4.41%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "114.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.31%  "
$ws.Range("B50").Value = "BitgetToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.67%  "
$ws.Range("D51").Value = "0.0₃0479"
$ws.Range("E51").Value = "  -9.32%  "
